$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.02743666666666666
$ws.Range("H2").Value = 0.08230999999999999
$ws.Range("I2").Value = 0.007366285056527356
$ws.Range("J2").Value = 0.007366285056527356
$ws.Range("M2").Value = 52.91852733333334
$ws.Range("N2").Value = 158.755582
$ws.Range("O2").Value = 0.9912603569328422
$ws.Range("P2").Value = 0.9912603569328421
$ws.Range("Q2").Value = 1.451907994935556
$ws.Range("R2").Value = 13.06717195442
$ws.Range("S2").Value = 0.007301906354402369
$ws.Range("T2").Value = 0.007301906354402368

$ws.Range("G3").Value = 0.02743666666666666
$ws.Range("H3").Value = 0.08230999999999999
$ws.Range("I3").Value = 0.007366285056527356
$ws.Range("J3").Value = 0.007366285056527356
$ws.Range("O3").Value = 0.003851187374513192
$ws.Range("P3").Value = 0.003851187374513192
$ws.Range("Q3").Value = 0.00564086892
$ws.Range("R3").Value = 0.05076782027999999
$ws.Range("S3").Value = 0.00002836894400676335
$ws.Range("T3").Value = 0.00002836894400676335

$ws.Range("G4").Value = 0.02743666666666666
$ws.Range("H4").Value = 0.08230999999999999
$ws.Range("I4").Value = 0.007366285056527356
$ws.Range("J4").Value = 0.007366285056527356
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.2609706666666667
$ws.Range("N4").Value = 0.7829120000000001
$ws.Range("O4").Value = 0.004888455692644593
$ws.Range("P4").Value = 0.004888455692644592
$ws.Range("Q4").Value = 0.007160165191111111
$ws.Range("R4").Value = 0.06444148672
$ws.Range("S4").Value = 0.00003600975811822395
$ws.Range("T4").Value = 0.00003600975811822394

$ws.Range("G5").Value = 3.368329
$ws.Range("I5").Value = 0.9043398704228307
$ws.Range("J5").Value = 0.9043398704228307
$ws.Range("M5").Value = 52.91852733333334
$ws.Range("N5").Value = 158.755582
$ws.Range("O5").Value = 0.9912603569328422
$ws.Range("P5").Value = 0.9912603569328421
$ws.Range("Q5").Value = 178.2470102541593
$ws.Range("R5").Value = 1604.223092287434
$ws.Range("S5").Value = 0.8964362627439354
$ws.Range("T5").Value = 0.8964362627439353

$ws.Range("G6").Value = 3.368329
$ws.Range("I6").Value = 0.9043398704228307
$ws.Range("J6").Value = 0.9043398704228307
$ws.Range("O6").Value = 0.003851187374513192
$ws.Range("P6").Value = 0.003851187374513192
$ws.Range("Q6").Value = 0.6925149690839999
$ws.Range("R6").Value = 6.232634721756
$ws.Range("S6").Value = 0.003482782291241302
$ws.Range("T6").Value = 0.003482782291241302

$ws.Range("G7").Value = 3.368329
$ws.Range("I7").Value = 0.9043398704228307
$ws.Range("J7").Value = 0.9043398704228307
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.2609706666666667
$ws.Range("N7").Value = 0.7829120000000001
$ws.Range("O7").Value = 0.004888455692644593
$ws.Range("P7").Value = 0.004888455692644592
$ws.Range("Q7").Value = 0.8790350646826667
$ws.Range("R7").Value = 7.911315582144
$ws.Range("S7").Value = 0.00442082538765396
$ws.Range("T7").Value = 0.004420825387653959

$ws.Range("G8").Value = 0.3288616666666667
$ws.Range("H8").Value = 0.9865849999999999
$ws.Range("I8").Value = 0.08829384452064198
$ws.Range("J8").Value = 0.08829384452064198
$ws.Range("M8").Value = 52.91852733333334
$ws.Range("N8").Value = 158.755582
$ws.Range("O8").Value = 0.9912603569328422
$ws.Range("P8").Value = 0.9912603569328421
$ws.Range("Q8").Value = 17.40287509638556
$ws.Range("R8").Value = 156.62587586747
$ws.Range("S8").Value = 0.08752218783450445
$ws.Range("T8").Value = 0.08752218783450444

$ws.Range("G9").Value = 0.3288616666666667
$ws.Range("H9").Value = 0.9865849999999999
$ws.Range("I9").Value = 0.08829384452064198
$ws.Range("J9").Value = 0.08829384452064198
$ws.Range("O9").Value = 0.003851187374513192
$ws.Range("P9").Value = 0.003851187374513192
$ws.Range("Q9").Value = 0.06761264322
$ws.Range("R9").Value = 0.60851378898
$ws.Range("S9").Value = 0.0003400361392651272
$ws.Range("T9").Value = 0.0003400361392651272

$ws.Range("G10").Value = 0.3288616666666667
$ws.Range("H10").Value = 0.9865849999999999
$ws.Range("I10").Value = 0.08829384452064198
$ws.Range("J10").Value = 0.08829384452064198
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.2609706666666667
$ws.Range("N10").Value = 0.7829120000000001
$ws.Range("O10").Value = 0.004888455692644593
$ws.Range("P10").Value = 0.004888455692644592
$ws.Range("Q10").Value = 0.08582324839111112
$ws.Range("R10").Value = 0.77240923552
$ws.Range("S10").Value = 0.0004316205468724089
$ws.Range("T10").Value = 0.0004316205468724088
